# "edit in upload functions"
# Updates the attendance/roster figures in column A (ID numbers) and the
# "annual" count in H4, relocates the active selection to H7, and sets the
# sheet's page setup (paper size / orientation) as it appears after the
# edit was made and the file re-saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated roster ID numbers (column A) ---------------------------------
$ws.Range("A2").Value = 1152
$ws.Range("A3").Value = 167
$ws.Range("A4").Value = 900

# --- Updated "annual" (H) leave count for row 4 ---------------------------
$ws.Range("H4").Value = 3

# --- Page setup: A4 paper, portrait orientation ----------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Move the active cell / selection to H7 --------------------------------
$ws.Range("H7").Select()
